$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed Price (D) / Volume(1h) (E) text values from the scheduled scraper run.
# These columns store plain text (e.g. "261.15", "-0.05%"), so the whole D2:E50 block is
# temporarily switched to Text format before writing, then reverted to General afterwards —
# otherwise Excel auto-coerces the numeric-looking / percent-looking strings into real numbers.
$dataRange = $ws.Range("D2:E50")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "261.15"
$ws.Range("E2").Value = "-0.05%"
$ws.Range("D3").Value = "26.91"
$ws.Range("E3").Value = "-1.49%"
$ws.Range("D4").Value = "4.721"
$ws.Range("E4").Value = "0.32%"
$ws.Range("D5").Value = "0.06206"
$ws.Range("E5").Value = "2.02%"
$ws.Range("D6").Value = "6.730"
$ws.Range("E6").Value = "0.87%"
$ws.Range("D7").Value = "0.8495"
$ws.Range("E7").Value = "0.41%"
$ws.Range("D8").Value = "0.9125"
$ws.Range("E8").Value = "-1.49%"
$ws.Range("D9").Value = "0.1403"
$ws.Range("E9").Value = "-0.01%"
$ws.Range("D10").Value = "0.04958"
$ws.Range("E10").Value = "4.89%"
$ws.Range("D11").Value = "0.07102"
$ws.Range("E11").Value = "-0.02%"
$ws.Range("D12").Value = "0.03118"
$ws.Range("E12").Value = "1.01%"
$ws.Range("D13").Value = "0.09052"
$ws.Range("E13").Value = "-0.16%"
$ws.Range("D14").Value = "0.001535"
$ws.Range("E14").Value = "-0.35%"
$ws.Range("D15").Value = "0.0006158"
$ws.Range("E15").Value = "1.19%"
$ws.Range("D16").Value = "0.005953"
$ws.Range("E16").Value = "-3.19%"
$ws.Range("E17").Value = "-0.05%"
$ws.Range("D18").Value = "3.174"
$ws.Range("E20").Value = "-0.37%"
$ws.Range("D21").Value = "0.1310"
$ws.Range("E21").Value = "1.67%"
$ws.Range("D22").Value = "4.091"
$ws.Range("E22").Value = "-0.18%"
$ws.Range("D23").Value = "0.04251"
$ws.Range("E23").Value = "0.13%"
$ws.Range("D24").Value = "0.001182"
$ws.Range("E24").Value = "-3.27%"
$ws.Range("D25").Value = "0.004066"
$ws.Range("E25").Value = "3.94%"
$ws.Range("E26").Value = "0.05%"
$ws.Range("E27").Value = "4.09%"
$ws.Range("D40").Value = "0.03944"
$ws.Range("E40").Value = "1.81%"
$ws.Range("E41").Value = "-0.12%"
$ws.Range("D42").Value = "0.004150"
$ws.Range("E42").Value = "0.96%"
$ws.Range("D43").Value = "0.002143"
$ws.Range("E43").Value = "-3.35%"
$ws.Range("E44").Value = "-19.50%"
$ws.Range("D45").Value = "0.00005165"
$ws.Range("E45").Value = "0.52%"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").Value = "0.05%"
$ws.Range("D48").Value = "0.2510"
$ws.Range("E48").Value = "85.50%"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").Value = "0.05%"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").Value = "0.05%"

$dataRange.ClearFormats()

Write-Host "Updated 63 cells."
